$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.183.91"
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").Value = "2.056.13"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.26"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.84"
$ws.Range("E7").Value = "  +7.88%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0807"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.81"
$ws.Range("E12").Value = "  +3.28%  "
$ws.Range("D13").Value = "2.360.43"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.23"
$ws.Range("E14").Value = "  +5.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.36"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.759"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").Value = "2.056.11"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "38.112.57"
$ws.Range("E18").Value = "  +2.50%  "
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").Value = ("0.0" + [string][char]0x2083 + "0833")
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.54"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.86"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("E28").Value = "  +3.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.93"
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("E31").Value = "  +2.17%  "
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("E34").Value = "  +4.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0606"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.44"
$ws.Range("E36").Value = "  +16.93%  "
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "1.518.89"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.59"
$ws.Range("E41").Value = "  +3.27%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.00"
$ws.Range("E42").Value = "  +4.13%  "
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.87"
$ws.Range("E44").Value = "  +3.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0926"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("E47").Value = "  -3.92%  "
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("D51").Value = "2.247.86"
$ws.Range("E51").Value = "  +1.54%  "
